# Update gh-pages to output generated at 456a3b4
# Applies the same updates to both the "展览" sheet and the "全部类型" sheet,
# since they contain duplicated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 3: 想去人数 (F3) 1395 -> 1400, Cover (I3) image URL updated
    $ws.Range("F3").Value = 1400
    $ws.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202410/UIoJJckw1729583568078.jpeg"

    # Row 4: 想去人数 (F4) 90 -> 92
    $ws.Range("F4").Value = 92
}
